$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Opportunities")
$ws.Range("A1").Value = "TC_ID"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = 65535
$ws.Range("A1").Borders.LineStyle = 1
